$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1399.75
$ws.Range("I48").Value = 1649.5
$ws.Range("J48").Value = 1150
$ws.Range("K48").Value = 4948.5
$ws.Range("L48").Value = 3450
$ws.Range("M48").Value = -4656.5
$ws.Range("N48").Value = -4034

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 285.77777
$ws.Range("I53").Value = 175.6
$ws.Range("J53").Value = 423.5
$ws.Range("K53").Value = 175.6
$ws.Range("L53").Value = 423.5
$ws.Range("M53").Value = 461.4
$ws.Range("N53").Value = -1697.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 1399.75
$ws.Range("I56").Value = 1649.5
$ws.Range("J56").Value = 1150
$ws.Range("K56").Value = 4948.5
$ws.Range("L56").Value = 3450
$ws.Range("M56").Value = -4414.5
$ws.Range("N56").Value = -4518

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 54372.2
$ws.Range("J70").Value = 75983.86
$ws.Range("L70").Value = 227951.58
$ws.Range("N70").Value = -228491.58

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 54372.2
$ws.Range("J73").Value = 75983.86
$ws.Range("L73").Value = 227951.58
$ws.Range("N73").Value = -229823.58

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3336.5574
$ws.Range("I138").Value = 1179.6666
$ws.Range("J138").Value = 5049.3823
$ws.Range("K138").Value = 3538.9998
$ws.Range("L138").Value = 15148.1469
$ws.Range("M138").Value = 1601.0002
$ws.Range("N138").Value = -25428.1469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10333.704
$ws.Range("I45").Value = 14156.235
$ws.Range("J45").Value = 3835.4
$ws.Range("K45").Value = 14156.235
$ws.Range("L45").Value = 3835.4
$ws.Range("M45").Value = -13779.235
$ws.Range("N45").Value = -4589.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 23888.334
$ws.Range("I53").Value = 18666.2
$ws.Range("K53").Value = 18666.2
$ws.Range("M53").Value = -17984.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5716.364
$ws.Range("I61").Value = 5057.706
$ws.Range("J61").Value = 7955.8
$ws.Range("K61").Value = 5057.706
$ws.Range("L61").Value = 7955.8
$ws.Range("M61").Value = -4845.706
$ws.Range("N61").Value = -8379.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3686
$ws.Range("I122").Value = 3583.25
$ws.Range("J122").Value = 3960
$ws.Range("K122").Value = 10749.75
$ws.Range("L122").Value = 11880
$ws.Range("M122").Value = -8299.75
$ws.Range("N122").Value = -16780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5716.364
$ws.Range("I136").Value = 5057.706
$ws.Range("J136").Value = 7955.8
$ws.Range("K136").Value = 15173.118
$ws.Range("L136").Value = 23867.4
$ws.Range("M136").Value = -12623.118
$ws.Range("N136").Value = -28967.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3997.2632
$ws.Range("I86").Value = 3144.2144
$ws.Range("J86").Value = 6385.8
$ws.Range("K86").Value = 3144.2144
$ws.Range("L86").Value = 6385.8
$ws.Range("M86").Value = -2021.2144
$ws.Range("N86").Value = -8631.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3997.2632
$ws.Range("I89").Value = 3144.2144
$ws.Range("J89").Value = 6385.8
$ws.Range("K89").Value = 15721.072
$ws.Range("L89").Value = 31929
$ws.Range("M89").Value = -10105.072
$ws.Range("N89").Value = -43161

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 89999
$ws.Range("J126").Value = 89999
$ws.Range("L126").Value = 89999
$ws.Range("N126").Value = -99879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2842.2144
$ws.Range("I134").Value = 3024.9546
$ws.Range("J134").Value = 2172.1667
$ws.Range("K134").Value = 9074.863799999999
$ws.Range("L134").Value = 6516.500100000001
$ws.Range("M134").Value = -6539.863799999999
$ws.Range("N134").Value = -11586.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9082.333000000001
$ws.Range("I86").Value = 8523.5
$ws.Range("K86").Value = 8523.5
$ws.Range("M86").Value = -7400.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 9082.333000000001
$ws.Range("I89").Value = 8523.5
$ws.Range("K89").Value = 42617.5
$ws.Range("M89").Value = -37001.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2744.6
$ws.Range("I107").Value = 2901.2307
$ws.Range("J107").Value = 1726.5
$ws.Range("K107").Value = 2901.2307
$ws.Range("L107").Value = 1726.5
$ws.Range("M107").Value = -981.2307000000001
$ws.Range("N107").Value = -5566.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 992.7
$ws.Range("I2").Value = 1437.091
$ws.Range("J2").Value = 449.55554
$ws.Range("K2").Value = 8622.545999999998
$ws.Range("L2").Value = 2697.33324
$ws.Range("M2").Value = -8509.545999999998
$ws.Range("N2").Value = -2923.33324

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 77284140
$ws.Range("J4").Value = 1334000.4
$ws.Range("L4").Value = 4002001.2
$ws.Range("N4").Value = -4002225.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 18000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 18000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 54000
$ws.Range("N9").Value = -54448
$ws.Range("M9").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1450
$ws.Range("J32").Value = 1400
$ws.Range("L32").Value = 4200
$ws.Range("N32").Value = -4766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 12498.5
$ws.Range("J100").Value = 12498.5
$ws.Range("L100").Value = 37495.5
$ws.Range("N100").Value = -39117.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 6402.4243
$ws.Range("I112").Value = 1570.5
$ws.Range("K112").Value = 4711.5
$ws.Range("M112").Value = -3603.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1876.5454
$ws.Range("I113").Value = 1950
$ws.Range("K113").Value = 5850
$ws.Range("M113").Value = -3680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 6941.5
$ws.Range("I126").Value = 2329.8
$ws.Range("K126").Value = 6989.400000000001
$ws.Range("M126").Value = -2049.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2062.4211
$ws.Range("I131").Value = 1245.3846
$ws.Range("K131").Value = 3736.1538
$ws.Range("M131").Value = 1303.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1483.1428
$ws.Range("I122").Value = 1483.1428
$ws.Range("K122").Value = 4449.428400000001
$ws.Range("M122").Value = -1999.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5377.6
$ws.Range("I126").Value = 5458.2856
$ws.Range("K126").Value = 16374.8568
$ws.Range("M126").Value = -13904.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12233.042
$ws.Range("I40").Value = 8174.067
$ws.Range("K40").Value = 8174.067
$ws.Range("M40").Value = -8038.067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 30046
$ws.Range("J48").Value = 30046
$ws.Range("L48").Value = 30046
$ws.Range("N48").Value = -31368

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4445
$ws.Range("I68").Value = 2621.75
$ws.Range("J68").Value = 5486.857
$ws.Range("K68").Value = 2621.75
$ws.Range("L68").Value = 5486.857
$ws.Range("M68").Value = -1872.75
$ws.Range("N68").Value = -6984.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4445
$ws.Range("I71").Value = 2621.75
$ws.Range("J71").Value = 5486.857
$ws.Range("K71").Value = 13108.75
$ws.Range("L71").Value = 27434.285
$ws.Range("M71").Value = -9364.75
$ws.Range("N71").Value = -34922.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1545.5
$ws.Range("I82").Value = 1374.5
$ws.Range("J82").Value = 1579.7
$ws.Range("K82").Value = 1374.5
$ws.Range("L82").Value = 1579.7
$ws.Range("M82").Value = -1013.5
$ws.Range("N82").Value = -2301.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1545.5
$ws.Range("I85").Value = 1374.5
$ws.Range("J85").Value = 1579.7
$ws.Range("K85").Value = 1374.5
$ws.Range("L85").Value = 1579.7
$ws.Range("M85").Value = -126.5
$ws.Range("N85").Value = -4075.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 48000
$ws.Range("J27").Value = 48000
$ws.Range("L27").Value = 48000
$ws.Range("N27").Value = -48138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18842.75
$ws.Range("J62").Value = 24915.666
$ws.Range("L62").Value = 24915.666
$ws.Range("N62").Value = -26163.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 18842.75
$ws.Range("J65").Value = 24915.666
$ws.Range("L65").Value = 124578.33
$ws.Range("N65").Value = -130818.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2215.9412
$ws.Range("I81").Value = 2289.182
$ws.Range("J81").Value = 2081.6667
$ws.Range("K81").Value = 4578.364
$ws.Range("L81").Value = 4163.3334
$ws.Range("M81").Value = -3517.364
$ws.Range("N81").Value = -6285.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2215.9412
$ws.Range("I84").Value = 2289.182
$ws.Range("J84").Value = 2081.6667
$ws.Range("K84").Value = 22891.82
$ws.Range("L84").Value = 20816.667
$ws.Range("M84").Value = -17587.82
$ws.Range("N84").Value = -31424.667
